$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44699
$ws.Range("M2").Value = 56
$ws.Range("N2").Value = 12000
$ws.Range("O2").Value = 12000
$ws.Range("P2").Value = 12000
$ws.Range("S2").Value = 1200
$ws.Range("D3").Value = 44699
$ws.Range("M3").Value = 60
$ws.Range("N3").Value = 10000
$ws.Range("O3").Value = 10000
$ws.Range("P3").Value = 10000
$ws.Range("S3").Value = 1000
$ws.Range("D4").Value = 44321
$ws.Range("L4").Value = 'Primera'
$ws.Range("M4").Value = 58
$ws.Range("N4").Value = 9000
$ws.Range("O4").Value = 9000
$ws.Range("P4").Value = 9000
$ws.Range("S4").Value = 900
$ws.Range("D5").Value = 44312
$ws.Range("M5").Value = 48
$ws.Range("D6").Value = 44319
$ws.Range("L6").Value = 'Primera'
$ws.Range("M6").Value = 68
$ws.Range("N6").Value = 10000
$ws.Range("O6").Value = 10000
$ws.Range("P6").Value = 10000
$ws.Range("S6").Value = 1000
$ws.Range("D7").Value = 44319
$ws.Range("L7").Value = 'Segunda'
$ws.Range("M7").Value = 57
$ws.Range("N7").Value = 8000
$ws.Range("O7").Value = 8000
$ws.Range("P7").Value = 8000
$ws.Range("S7").Value = 800
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 65
$ws.Range("N8").Value = 10000
$ws.Range("O8").Value = 10000
$ws.Range("P8").Value = 10000
$ws.Range("S8").Value = 1000
$ws.Range("D9").Value = 44326
$ws.Range("L9").Value = 'Segunda'
$ws.Range("M9").Value = 67
$ws.Range("N9").Value = 8000
$ws.Range("O9").Value = 8000
$ws.Range("P9").Value = 8000
$ws.Range("S9").Value = 800
$ws.Range("D10").Value = 44323
$ws.Range("M10").Value = 60
$ws.Range("D11").Value = 44323
$ws.Range("M11").Value = 50
$ws.Range("N11").Value = 9000
$ws.Range("O11").Value = 9000
$ws.Range("P11").Value = 9000
$ws.Range("S11").Value = 900
$ws.Range("D12").Value = 44301
$ws.Range("M12").Value = 45
$ws.Range("D13").Value = 44308
$ws.Range("L13").Value = 'Primera'
$ws.Range("M13").Value = 45
$ws.Range("N13").Value = 10000
$ws.Range("O13").Value = 10000
$ws.Range("P13").Value = 10000
$ws.Range("S13").Value = 1000
$ws.Range("D14").Value = 44308
$ws.Range("L14").Value = 'Segunda'
$ws.Range("M14").Value = 48
$ws.Range("D15").Value = 44315
$ws.Range("L15").Value = 'Primera'
$ws.Range("M15").Value = 45
$ws.Range("N15").Value = 10000
$ws.Range("O15").Value = 10000
$ws.Range("P15").Value = 10000
$ws.Range("S15").Value = 1000
$ws.Range("D16").Value = 44314
$ws.Range("L16").Value = 'Primera'
$ws.Range("M16").Value = 47
$ws.Range("N16").Value = 9000
$ws.Range("O16").Value = 9000
$ws.Range("P16").Value = 9000
$ws.Range("S16").Value = 900
$ws.Range("D17").Value = 44306
$ws.Range("M17").Value = 45
$ws.Range("D21").Value = 44309
$ws.Range("M21").Value = 45
$ws.Range("D22").Value = 44307
$ws.Range("M22").Value = 40
$ws.Range("D23").Value = 44329
$ws.Range("M23").Value = 56
$ws.Range("N23").Value = 9000
$ws.Range("O23").Value = 9000
$ws.Range("P23").Value = 9000
$ws.Range("R23").Value = 'Región Metropolitana'
$ws.Range("S23").Value = 900
$ws.Range("D24").Value = 44329
$ws.Range("M24").Value = 50
$ws.Range("R24").Value = 'Región Metropolitana'
$ws.Range("D25").Value = 44328
$ws.Range("M25").Value = 45
$ws.Range("N25").Value = 8000
$ws.Range("O25").Value = 8000
$ws.Range("P25").Value = 8000
$ws.Range("S25").Value = 800
$ws.Range("D26").Value = 44328
$ws.Range("L26").Value = 'Segunda'
$ws.Range("M26").Value = 48
$ws.Range("N26").Value = 7000
$ws.Range("O26").Value = 7000
$ws.Range("P26").Value = 7000
$ws.Range("S26").Value = 700
$ws.Range("D27").Value = 44333
$ws.Range("L27").Value = 'Especial'
$ws.Range("M27").Value = 58
$ws.Range("D28").Value = 44333
$ws.Range("M28").Value = 65
$ws.Range("R28").Value = 'Provincia de Quillota'
$ws.Range("D29").Value = 44333
$ws.Range("M29").Value = 60
$ws.Range("R29").Value = 'Provincia de Quillota'
$ws.Range("D30").Value = 44322
$ws.Range("M30").Value = 56
$ws.Range("D31").Value = 44322
$ws.Range("L31").Value = 'Segunda'
$ws.Range("M31").Value = 40
$ws.Range("N31").Value = 8000
$ws.Range("O31").Value = 8000
$ws.Range("P31").Value = 8000
$ws.Range("S31").Value = 800
$ws.Range("D32").Value = 44302
